$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y, Z, AA, AB hold date/time-like text (e.g. "2012-06-25", "00:00").
# Force them to Text format first so Excel does not auto-convert the literal
# strings we are about to write into date/time serial numbers.
$ws.Range("Y2:AB8").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 2682473
$ws.Range("B2").Value = 101691
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 220164
$ws.Range("F2").Value = "Solvända"
$ws.Range("G2").Value = "Helianthemum nummularium"
$ws.Range("H2").Value = "(L.) Mill."
$ws.Range("P2").Value = "150 m SSÖ om Örsta, vägkanten i kurvan, Upl"
$ws.Range("Q2").Value = 676571.1877922384
$ws.Range("R2").Value = 6617992.858972132
$ws.Range("S2").Value = 25
$ws.Range("T2").Value = "Stockholm"
$ws.Range("U2").Value = "Vallentuna"
$ws.Range("V2").Value = "Uppland"
$ws.Range("W2").Value = "Frösunda"
$ws.Range("Y2").Value = "2012-06-25"
$ws.Range("Z2").Value = "00:00"
$ws.Range("AA2").Value = "2012-06-25"
$ws.Range("AB2").Value = "00:00"
$ws.Range("AD2").Value = $False
$ws.Range("AE2").Value = $False
$ws.Range("AG2").Value = $False
$ws.Range("AW2").Value = "Måns Svensson"
$ws.Range("AX2").Value = "Måns Svensson"
# Row 3
$ws.Range("A3").Value = 4197334
$ws.Range("B3").Value = 104404
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 221849
$ws.Range("F3").Value = "Backtimjan"
$ws.Range("G3").Value = "Thymus serpyllum"
$ws.Range("H3").Value = "L."
$ws.Range("P3").Value = "150 m SSÖ om Örsta, vägkanten i kurvan, Upl"
$ws.Range("Q3").Value = 676571.1877922384
$ws.Range("R3").Value = 6617992.858972132
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = "Stockholm"
$ws.Range("U3").Value = "Vallentuna"
$ws.Range("V3").Value = "Uppland"
$ws.Range("W3").Value = "Frösunda"
$ws.Range("Y3").Value = "2012-06-25"
$ws.Range("Z3").Value = "00:00"
$ws.Range("AA3").Value = "2012-06-25"
$ws.Range("AB3").Value = "00:00"
$ws.Range("AD3").Value = $False
$ws.Range("AE3").Value = $False
$ws.Range("AG3").Value = $False
$ws.Range("AW3").Value = "Måns Svensson"
$ws.Range("AX3").Value = "Måns Svensson"
# Row 4
$ws.Range("A4").Value = 111634290
$ws.Range("B4").Value = 98535
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = "Blåsippa"
$ws.Range("G4").Value = "Hepatica nobilis"
$ws.Range("H4").Value = "Schreb."
$ws.Range("P4").Value = "Örsta, Upl"
$ws.Range("Q4").Value = 676709
$ws.Range("R4").Value = 6618511
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Stockholm"
$ws.Range("U4").Value = "Vallentuna"
$ws.Range("V4").Value = "Uppland"
$ws.Range("W4").Value = "Frösunda"
$ws.Range("Y4").Value = "2023-08-22"
$ws.Range("Z4").Value = "17:43"
$ws.Range("AA4").Value = "2023-08-22"
$ws.Range("AB4").Value = "17:43"
$ws.Range("AD4").Value = $False
$ws.Range("AE4").Value = $False
$ws.Range("AG4").Value = $False
$ws.Range("AW4").Value = "Karolin Hård"
$ws.Range("AX4").Value = "Karolin Hård"
# Row 5
$ws.Range("A5").Value = 111633890
$ws.Range("B5").Value = 90658
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 4361
$ws.Range("F5").Value = "Orange taggsvamp"
$ws.Range("G5").Value = "Hydnellum aurantiacum"
$ws.Range("H5").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P5").Value = "Örsta, Upl"
$ws.Range("Q5").Value = 676487
$ws.Range("R5").Value = 6618440
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Stockholm"
$ws.Range("U5").Value = "Vallentuna"
$ws.Range("V5").Value = "Uppland"
$ws.Range("W5").Value = "Frösunda"
$ws.Range("Y5").Value = "2023-08-22"
$ws.Range("Z5").Value = "17:43"
$ws.Range("AA5").Value = "2023-08-22"
$ws.Range("AB5").Value = "17:43"
$ws.Range("AD5").Value = $False
$ws.Range("AE5").Value = $False
$ws.Range("AG5").Value = $False
$ws.Range("AW5").Value = "Karolin Hård"
$ws.Range("AX5").Value = "Karolin Hård"
# Row 6
$ws.Range("A6").Value = 111633843
$ws.Range("B6").Value = 90687
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 5964
$ws.Range("F6").Value = "Fjällig taggsvamp s.str."
$ws.Range("G6").Value = "Sarcodon imbricatus s.str."
$ws.Range("H6").Value = "(L.:Fr.) P.Karst."
$ws.Range("P6").Value = "Örsta, Upl"
$ws.Range("Q6").Value = 676487
$ws.Range("R6").Value = 6618440
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Stockholm"
$ws.Range("U6").Value = "Vallentuna"
$ws.Range("V6").Value = "Uppland"
$ws.Range("W6").Value = "Frösunda"
$ws.Range("Y6").Value = "2023-08-22"
$ws.Range("Z6").Value = "17:43"
$ws.Range("AA6").Value = "2023-08-22"
$ws.Range("AB6").Value = "17:43"
$ws.Range("AD6").Value = $False
$ws.Range("AE6").Value = $False
$ws.Range("AG6").Value = $False
$ws.Range("AW6").Value = "Karolin Hård"
$ws.Range("AX6").Value = "Karolin Hård"
# Row 7
$ws.Range("A7").Value = 111633837
$ws.Range("B7").Value = 98535
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 222498
$ws.Range("F7").Value = "Blåsippa"
$ws.Range("G7").Value = "Hepatica nobilis"
$ws.Range("H7").Value = "Schreb."
$ws.Range("P7").Value = "Örsta, Upl"
$ws.Range("Q7").Value = 676487
$ws.Range("R7").Value = 6618440
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = "Stockholm"
$ws.Range("U7").Value = "Vallentuna"
$ws.Range("V7").Value = "Uppland"
$ws.Range("W7").Value = "Frösunda"
$ws.Range("Y7").Value = "2023-08-22"
$ws.Range("Z7").Value = "17:43"
$ws.Range("AA7").Value = "2023-08-22"
$ws.Range("AB7").Value = "17:43"
$ws.Range("AD7").Value = $False
$ws.Range("AE7").Value = $False
$ws.Range("AG7").Value = $False
$ws.Range("AW7").Value = "Karolin Hård"
$ws.Range("AX7").Value = "Karolin Hård"
# Row 8
$ws.Range("A8").Value = 111634304
$ws.Range("B8").Value = 90687
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 5964
$ws.Range("F8").Value = "Fjällig taggsvamp s.str."
$ws.Range("G8").Value = "Sarcodon imbricatus s.str."
$ws.Range("H8").Value = "(L.:Fr.) P.Karst."
$ws.Range("P8").Value = "Örsta, Upl"
$ws.Range("Q8").Value = 676709
$ws.Range("R8").Value = 6618511
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = "Stockholm"
$ws.Range("U8").Value = "Vallentuna"
$ws.Range("V8").Value = "Uppland"
$ws.Range("W8").Value = "Frösunda"
$ws.Range("Y8").Value = "2023-08-22"
$ws.Range("Z8").Value = "17:43"
$ws.Range("AA8").Value = "2023-08-22"
$ws.Range("AB8").Value = "17:43"
$ws.Range("AD8").Value = $False
$ws.Range("AE8").Value = $False
$ws.Range("AG8").Value = $False
$ws.Range("AW8").Value = "Karolin Hård"
$ws.Range("AX8").Value = "Karolin Hård"

# Column K ("Ålder-Stadium") previously had empty placeholder cells on rows
# 2-6 and none on rows 7-8. After the edit, rows 2 and 3 should no longer
# carry that empty cell, while rows 7 and 8 should gain one (matching the
# rotated row contents).
$ws.Range("K2").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("K4").Copy($ws.Range("K7"))
$ws.Range("K4").Copy($ws.Range("K8"))
